# Auto-generated script to update cryptos.xlsx per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.779.68'  # D2: was '42.831.20'
$ws.Cells.Item(2, 5).Value = '  +0.43%  '  # E2: was '  +0.57%  '
$ws.Cells.Item(3, 4).Value = '2.547.91'  # D3: was '2.542.65'
$ws.Cells.Item(3, 5).Value = '  +0.13%  '  # E3: was '  -0.09%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '  # E4: was '  -0.15%  '
$ws.Cells.Item(5, 4).Value = '''304.23'  # D5: was '303.78'
$ws.Cells.Item(5, 5).Value = '  +1.55%  '  # E5: was '  +1.60%  '
$ws.Cells.Item(6, 4).Value = '''97.83'  # D6: was '97.78'
$ws.Cells.Item(6, 5).Value = '  +4.25%  '  # E6: was '  +6.12%  '
$ws.Cells.Item(7, 5).Value = '  +0.42%  '  # E7: was '  +0.41%  '
$ws.Cells.Item(8, 5).Value = '  +0.07%  '  # E8: was '  +0.04%  '
$ws.Cells.Item(9, 4).Value = '''0.546'  # D9: was '0.545'
$ws.Cells.Item(9, 5).Value = '  -0.20%  '  # E9: was '  -0.76%  '
$ws.Cells.Item(10, 4).Value = '''36.79'  # D10: was '36.84'
$ws.Cells.Item(10, 5).Value = '  +2.27%  '  # E10: was '  +2.69%  '
$ws.Cells.Item(11, 4).Value = '''0.0828'  # D11: was '0.0830'
$ws.Cells.Item(11, 5).Value = '  +3.03%  '  # E11: was '  +3.32%  '
$ws.Cells.Item(12, 4).Value = '''7.71'  # D12: was '7.72'
$ws.Cells.Item(12, 5).Value = '  -0.03%  '  # E12: was '  +0.77%  '
$ws.Cells.Item(13, 4).Value = '''0.115'  # D13: was '0.114'
$ws.Cells.Item(13, 5).Value = '  +1.60%  '  # E13: was '  +0.84%  '
$ws.Cells.Item(14, 4).Value = '2.936.60'  # D14: was '2.934.64'
$ws.Cells.Item(14, 5).Value = '  +0.20%  '  # E14: was '  +0.02%  '
$ws.Cells.Item(15, 4).Value = '2.555.76'  # D15: was '2.536.75'
$ws.Cells.Item(15, 5).Value = '  +0.08%  '  # E15: was '  -0.27%  '
$ws.Cells.Item(16, 4).Value = '''15.04'  # D16: was '15.09'
$ws.Cells.Item(16, 5).Value = '  +5.39%  '  # E16: was '  +5.98%  '
$ws.Cells.Item(17, 4).Value = '''0.867'  # D17: was '0.866'
$ws.Cells.Item(17, 5).Value = '  -0.54%  '  # E17: was '  -0.85%  '
$ws.Cells.Item(18, 4).Value = '42.827.23'  # D18: was '42.840.65'
$ws.Cells.Item(18, 5).Value = '  +0.48%  '  # E18: was '  +0.53%  '
$ws.Cells.Item(19, 4).Value = '''13.30'  # D19: was '13.34'
$ws.Cells.Item(19, 5).Value = '  +4.22%  '  # E19: was '  +3.80%  '
$ws.Cells.Item(20, 5).Value = '  +0.58%  '  # E20: was '  +0.84%  '
$ws.Cells.Item(21, 5).Value = '  +0.35%  '  # E21: was '  +0.49%  '
$ws.Cells.Item(22, 4).Value = '''71.82'  # D22: was '71.80'
$ws.Cells.Item(22, 5).Value = '  +0.11%  '  # E22: was '  +0.59%  '
$ws.Cells.Item(23, 4).Value = '''255.97'  # D23: was '255.98'
$ws.Cells.Item(23, 5).Value = '  -0.04%  '  # E23: was '  +0.13%  '
$ws.Cells.Item(24, 5).Value = '  +1.02%  '  # E24: was '  +1.28%  '
$ws.Cells.Item(25, 5).Value = '  -2.11%  '  # E25: was '  -1.62%  '
$ws.Cells.Item(26, 4).Value = '''28.05'  # D26: was '28.09'
$ws.Cells.Item(26, 5).Value = '  -3.85%  '  # E26: was '  -3.60%  '
$ws.Cells.Item(27, 5).Value = '  +0.09%  '  # E27: was '  -0.08%  '
$ws.Cells.Item(28, 4).Value = '''2.30'  # D28: was '2.32'
$ws.Cells.Item(28, 5).Value = '  +8.70%  '  # E28: was '  +9.26%  '
$ws.Cells.Item(29, 4).Value = '''10.17'  # D29: was '10.16'
$ws.Cells.Item(29, 5).Value = '  +1.33%  '  # E29: was '  +1.21%  '
$ws.Cells.Item(30, 4).Value = '''37.99'  # D30: was '37.93'
$ws.Cells.Item(30, 5).Value = '  +3.14%  '  # E30: was '  +2.84%  '
$ws.Cells.Item(31, 4).Value = '''6.12'  # D31: was '6.08'
$ws.Cells.Item(31, 5).Value = '  +3.05%  '  # E31: was '  +2.12%  '
$ws.Cells.Item(32, 4).Value = '''157.35'  # D32: was '157.71'
$ws.Cells.Item(32, 5).Value = '  +3.52%  '  # E32: was '  +3.25%  '
$ws.Cells.Item(33, 4).Value = '''19.70'  # D33: was '19.47'
$ws.Cells.Item(33, 5).Value = '  +15.24%  '  # E33: was '  +13.51%  '
$ws.Cells.Item(34, 4).Value = '''2.14'  # D34: was '2.13'
$ws.Cells.Item(34, 5).Value = '  -0.65%  '  # E34: was '  -1.48%  '
$ws.Cells.Item(35, 5).Value = '  +0.78%  '  # E35: was '  +0.77%  '
$ws.Cells.Item(36, 5).Value = '  -2.14%  '  # E36: was '  -2.54%  '
$ws.Cells.Item(37, 5).Value = '  -4.63%  '  # E37: was '  -4.24%  '
$ws.Cells.Item(38, 4).Value = '''25.51'  # D38: was '25.90'
$ws.Cells.Item(38, 5).Value = '  +4.43%  '  # E38: was '  +8.34%  '
$ws.Cells.Item(39, 5).Value = '  +1.30%  '  # E39: was '  +1.42%  '
$ws.Cells.Item(40, 5).Value = '  +0.43%  '  # E40: was '  +0.23%  '
$ws.Cells.Item(41, 4).Value = '''2.10'  # D41: was '2.09'
$ws.Cells.Item(41, 5).Value = '  +29.43%  '  # E41: was '  +29.34%  '
$ws.Cells.Item(42, 2).Value = 'NEARProtocol'  # B42: was 'RenderToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'  # C42: was 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(42, 4).Value = '''3.42'  # D42: was '3.89'
$ws.Cells.Item(42, 5).Value = '  +0.96%  '  # E42: was '  +0.59%  '
$ws.Cells.Item(43, 2).Value = 'RenderToken'  # B43: was 'NEARProtocol'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'  # C43: was 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(43, 4).Value = '''3.88'  # D43: was '3.40'
$ws.Cells.Item(43, 5).Value = '  +0.47%  '  # E43: was '  +0.24%  '
$ws.Cells.Item(44, 4).Value = '2.089.20'  # D44: was '2.091.36'
$ws.Cells.Item(44, 5).Value = '  +0.23%  '  # E44: was '  +0.20%  '
$ws.Cells.Item(45, 5).Value = '  -1.55%  '  # E45: was '  -1.44%  '
$ws.Cells.Item(46, 5).Value = '  +0.04%  '  # E46: was '  +0.02%  '
$ws.Cells.Item(47, 4).Value = '''87.16'  # D47: was '87.68'
$ws.Cells.Item(47, 5).Value = '  +3.11%  '  # E47: was '  +3.91%  '
$ws.Cells.Item(48, 4).Value = '''8.92'  # D48: was '8.86'
$ws.Cells.Item(48, 5).Value = '  -1.33%  '  # E48: was '  -3.24%  '
$ws.Cells.Item(49, 4).Value = '2.794.24'  # D49: was '2.792.34'
$ws.Cells.Item(49, 5).Value = '  +0.25%  '  # E49: was '  +0.06%  '
$ws.Cells.Item(50, 4).Value = '''74.52'  # D50: was '74.49'
$ws.Cells.Item(50, 5).Value = '  +8.04%  '  # E50: was '  +8.23%  '
$ws.Cells.Item(51, 2).Value = 'Aave'  # B51: was 'Algorand'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'  # C51: was 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(51, 4).Value = '''103.17'  # D51: was '0.191'
$ws.Cells.Item(51, 5).Value = '  -0.80%  '  # E51: was '  +2.71%  '
